$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.405.37'
$ws.Range("E2").Value = '  -2.69%  '

$ws.Range("D3").Value = '3.093.30'
$ws.Range("E3").Value = '  -1.44%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '547.63'
$ws.Range("E5").Value = '  -2.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.95'
$ws.Range("E6").Value = '  -6.01%  '

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").Value = '3.083.07'
$ws.Range("E8").Value = '  -1.52%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.494'
$ws.Range("E9").Value = '  -0.50%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.57'
$ws.Range("E10").Value = '  -4.88%  '

$ws.Range("E11").Value = '  +0.54%  '

$ws.Range("E12").Value = '  +0.40%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.01'
$ws.Range("E13").Value = '  -2.57%  '

$ws.Range("E14").Value = '  -0.82%  '

$ws.Range("D15").Value = '3.591.85'
$ws.Range("E15").Value = '  -1.33%  '

$ws.Range("D16").Value = '63.470.82'
$ws.Range("E16").Value = '  -2.52%  '

$ws.Range("E17").Value = '  -0.82%  '

$ws.Range("D18").Value = '3.091.39'
$ws.Range("E18").Value = '  -1.50%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '494.27'
$ws.Range("E19").Value = '  -5.00%  '

$ws.Range("E20").Value = '  -0.55%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.49'
$ws.Range("E21").Value = '  -2.05%  '

$ws.Range("E22").Value = '  +0.71%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.22'
$ws.Range("E23").Value = '  -2.52%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.93'
$ws.Range("E24").Value = '  -0.77%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.29'
$ws.Range("E25").Value = '  -2.66%  '

$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.74'
$ws.Range("E27").Value = '  -1.36%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.35'
$ws.Range("E28").Value = '  -2.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '26.56'
$ws.Range("E30").Value = '  +2.13%  '

$ws.Range("E31").Value = '  -8.23%  '

$ws.Range("E32").Value = '  +0.86%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '61.48'
$ws.Range("E33").Value = '  +16.27%  '

$ws.Range("E34").Value = '  -5.72%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '528.90'
$ws.Range("E35").Value = '  -5.97%  '

$ws.Range("E36").Value = '  -1.10%  '

$ws.Range("E37").Value = '  -4.36%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0405'
$ws.Range("E38").Value = '  -6.62%  '

$ws.Range("E39").Value = '  -2.45%  '

$ws.Range("D40").Value = '3.066.70'
$ws.Range("E40").Value = '  -0.35%  '

$ws.Range("E41").Value = '  -0.86%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.14'
$ws.Range("E42").Value = '  -0.51%  '

$ws.Range("E43").Value = '  -6.68%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.258'
$ws.Range("E44").Value = '  +0.77%  '

$ws.Range("E45").Value = '  +0.11%  '

$ws.Range("E46").Value = '  -5.95%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '121.98'
$ws.Range("E47").Value = '  +3.36%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.50'
$ws.Range("E48").Value = '  -1.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.108'
$ws.Range("E49").Value = '  +0.02%  '

$ws.Range("D50").Value = '0.0₃0509'
$ws.Range("E50").Value = '  -2.48%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.37'
$ws.Range("E51").Value = '  +54.40%  '
